$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (CI_Type, CI_Name, CI_Descrip, Rel_Type, Dependency_Type, Dependency_Name, Dependency_Descrip)
$newRows = @(
    @("Applications", "App 2", "Parent Description…", "Depends On", "Facilities", "Location 2", "Dependency Description…"),
    @("Applications", "App 2", "Parent Description…", "Depends On", "Procurements", "PO 2", "Dependency Description…"),
    @("Applications", "App 2", "Parent Description…", "Depends On", "People", "Person 2", "Dependency Description…"),
    @("Applications", "App 2", "Parent Description…", "Depends On", "Data", "Data 2", "Dependency Description…"),
    @("Applications", "App 2", "Parent Description…", "Depends On", "Technology", "Tech 2", "Dependency Description…")
)

$startRow = 32
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

# Widen column A to fit its content (bestFit-like behavior)
$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null

# Update the view: scroll so row 13 is the top-left visible row, and select E22
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("E22").Select() | Out-Null
